# Updated from BAG incidence numbers
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 4: demographics - refresh incidence (M) and replace the
# "cases per 100k" (N) column with a Canton / ": " / value layout
# (O / P / Q) using un-shared formulas.
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()

$ws4.Range("M2").Value = 30.1
$ws4.Range("M3").Value = 20.8
$ws4.Range("M4").Value = 115
$ws4.Range("M5").Value = 14.2
$ws4.Range("M6").Value = 17.3
$ws4.Range("M7").Value = 61.9
$ws4.Range("M8").Value = 15.4
$ws4.Range("M9").Value = 178.6
$ws4.Range("M10").Value = 51.2
$ws4.Range("M11").Value = 25.7
$ws4.Range("M12").Value = 48.2
$ws4.Range("M13").Value = 10.1
$ws4.Range("M14").Value = 14.3
$ws4.Range("M15").Value = 62
$ws4.Range("M16").Value = 122.7
$ws4.Range("M17").Value = 49.2
$ws4.Range("M18").Value = 22
$ws4.Range("M19").Value = 18.1
$ws4.Range("M20").Value = 9.8
$ws4.Range("M21").Value = 34.1
$ws4.Range("M22").Value = 29
$ws4.Range("M23").Value = 46.3
$ws4.Range("M24").Value = 39.6
$ws4.Range("M25").Value = 44.9
# M26 (GE) and M27 (JU) are unchanged

# Drop the old "cases per 100k" (N) column entirely
$ws4.Range("N2:N28").ClearContents()

# New "label" column: literal ": " shared string next to Canton (O)
$ws4.Range("P2:P28").Value = ": "

# New (non-shared) formula column replacing N, now in Q
$ws4.Range("Q2").Formula = "=ROUND(L2*M2,0)"
$ws4.Range("Q3").Formula = "=ROUND(L3*M3,0)"
$ws4.Range("Q4").Formula = "=ROUND(L4*M4,0)"
$ws4.Range("Q5").Formula = "=ROUND(L5*M5,0)"
$ws4.Range("Q6").Formula = "=ROUND(L6*M6,0)"
$ws4.Range("Q7").Formula = "=ROUND(L7*M7,0)"
$ws4.Range("Q8").Formula = "=ROUND(L8*M8,0)"
$ws4.Range("Q9").Formula = "=ROUND(L9*M9,0)"
$ws4.Range("Q10").Formula = "=ROUND(L10*M10,0)"
$ws4.Range("Q11").Formula = "=ROUND(L11*M11,0)"
$ws4.Range("Q12").Formula = "=ROUND(L12*M12,0)"
$ws4.Range("Q13").Formula = "=ROUND(L13*M13,0)"
$ws4.Range("Q14").Formula = "=ROUND(L14*M14,0)"
$ws4.Range("Q15").Formula = "=ROUND(L15*M15,0)"
$ws4.Range("Q16").Formula = "=ROUND(L16*M16,0)"
$ws4.Range("Q17").Formula = "=ROUND(L17*M17,0)"
$ws4.Range("Q18").Formula = "=ROUND(L18*M18,0)"
$ws4.Range("Q19").Formula = "=ROUND(L19*M19,0)"
$ws4.Range("Q20").Formula = "=ROUND(L20*M20,0)"
$ws4.Range("Q21").Formula = "=ROUND(L21*M21,0)"
$ws4.Range("Q22").Formula = "=ROUND(L22*M22,0)"
$ws4.Range("Q23").Formula = "=ROUND(L23*M23,0)"
$ws4.Range("Q24").Formula = "=ROUND(L24*M24,0)"
$ws4.Range("Q25").Formula = "=ROUND(L25*M25,0)"
$ws4.Range("Q26").Formula = "=ROUND(L26*M26,0)"
$ws4.Range("Q27").Formula = "=ROUND(L27*M27,0)"
$ws4.Range("Q28").Formula = "=SUM(Q2:Q27)"

# Restore the sheet selection/scroll like the author left it
$ws4.Range("O2:Q28").Select()

# -----------------------------------------------------------------
# Sheet 1: covid19_cases_switzerland - fill in row 15 (2020-03-18)
# (done last so this sheet ends up the active/selected one again,
# matching the original tabSelected="1" on this sheet)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

$ws1.Range("C15").Value = 3
$ws1.Range("D15").Value = 16
$ws1.Range("E15").Value = 215
$ws1.Range("F15").Value = 139
$ws1.Range("H15").Value = 82
$ws1.Range("J15").Value = 16
$ws1.Range("K15").Value = 123
$ws1.Range("L15").Value = 25
$ws1.Range("M15").Value = 63
$ws1.Range("O15").Value = 20
$ws1.Range("P15").Value = 17
$ws1.Range("R15").Value = 8
$ws1.Range("S15").Value = 39
$ws1.Range("T15").Value = 35
$ws1.Range("X15").Value = 919
$ws1.Range("Y15").Value = 176
$ws1.Range("Z15").Value = 23
# CH total for the day, was an empty (but number-formatted) cell
$ws1.Range("AB15").ClearFormats()
$ws1.Range("AB15").Value = 3855

# Move/restore the sheet selection like the author left it
$ws1.Range("A16").Select()
